# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '35.350.85'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +0.43%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.907.95'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +2.68%  '
$ws.Range('E4').Value = '  -0.43%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '245.75'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +2.81%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.663'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +6.51%  '
$ws.Range('E7').Value = '  -0.42%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '41.37'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -2.06%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.349'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +6.12%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '52.82'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +12.63%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0718'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +3.68%  '
$ws.Range('E12').Value = '  +0.40%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.186.71'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +2.78%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '12.09'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +5.14%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '1.924.93'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +3.29%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.699'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +3.31%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '4.86'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +2.88%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '35.338.20'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +0.50%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '72.58'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +3.90%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.0₃0828'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +4.16%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '239.66'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -0.42%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '12.55'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('E25').Value = '  +0.93%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +23.43%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '169.90'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +0.12%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '8.46'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +5.61%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '18.44'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +4.49%  '
$ws.Range('E30').Value = '  +2.28%  '
$ws.Range('E31').Value = '  +3.64%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.0565'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.939'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +15.25%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -0.28%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '4.11'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('E36').Value = '  -4.52%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('E40').Value = '  +3.32%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '16.30'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +9.09%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.0637'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +6.33%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '89.90'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -0.12%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.339.16'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -0.53%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.38'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +2.76%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '47.77'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +37.59%  '
$ws.Range('E47').Value = '  +1.80%  '
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('E49').Value = '  -0.37%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.092.79'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('E51').Value = '  +3.86%  '
